$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet2 (Second batch) cell updates ---
$ws2.Range('D2').Value = 'https://www.linkedin.com/in/gaurav-pandey-412180194/'
$ws2.Range('E2').Value = '../members/gaurav.jpg'
$ws2.Range('D3').Value = 'https://www.linkedin.com/in/kunal-thakur-b27823193/'
$ws2.Range('E3').Value = '../members/kunalt.jpg'
$ws2.Range('D4').Value = 'https://www.linkedin.com/in/ayushi-sharma-778657198/'
$ws2.Range('E4').Value = '../members/ayushi.jpg'
$ws2.Range('D5').Value = 'https://www.linkedin.com/in/jai-gupta-8b5238196/'
$ws2.Range('E5').Value = '../members/jai.jpg'
$ws2.Range('D6').Value = 'https://www.instagram.com/mnvndra/'
$ws2.Range('E6').Value = '../members/manvendra.jpg'
$ws2.Range('A7').Value = 'Mayur Kumar'
$ws2.Range('D7').Value = 'https://www.linkedin.com/in/mayur-kumar-47a9511a1/'
$ws2.Range('E7').Value = '../members/mayur.jpg'
$ws2.Range('D8').Value = 'https://www.linkedin.com/in/parasaggarwal/'
$ws2.Range('E8').Value = '../members/paras.jpg'
$ws2.Range('D9').Value = 'https://www.linkedin.com/in/parth-pant-866bb4189/'
$ws2.Range('E9').Value = '../members/parth.jpg'
$ws2.Range('D11').Value = 'https://www.linkedin.com/in/suryansh-dwivedi-9b9975199/'
$ws2.Range('E11').Value = '../members/suryansh.jpg'
$ws2.Range('D12').Value = 'https://www.linkedin.com/in/utkarsh-rai-50738b61/'
$ws2.Range('E12').Value = '../members/utkarsh.jpg'
$ws2.Range('D13').Value = 'https://www.linkedin.com/in/suryansh-dwivedi-9b9975199/'
$ws2.Range('D14').Value = 'https://www.linkedin.com/in/vanshika-thakur-6a7b4a19a/'
$ws2.Range('E14').Value = '../members/vanshika.jpg'
$ws2.Range('D15').Value = 'https://www.linkedin.com/in/vasundhra-thakur-5b9023196/'
$ws2.Range('E15').Value = '../members/vasundhra.jpg'

# --- Sheet3 (Third batch) cell updates ---
$ws3.Range('D2').Value = 'https://www.linkedin.com/in/yash-punia/'
$ws3.Range('E2').Value = '../members/yash.jpg'
$ws3.Range('D3').Value = 'https://www.linkedin.com/in/akhyarai/'
$ws3.Range('E3').Value = '../members/akhya.jpg'
$ws3.Range('D4').Value = 'https://www.linkedin.com/in/rathod-sahaab/'
$ws3.Range('E4').Value = '../members/abhay.jpg'
$ws3.Range('D5').Value = 'https://www.linkedin.com/in/aditi-singh2000/'
$ws3.Range('E5').Value = '../members/aditi.jpg'
$ws3.Range('D6').Value = 'https://www.linkedin.com/in/anshudhar-kumar-singh/'
$ws3.Range('E6').Value = '../members/anshudhar.jpg'
$ws3.Range('D7').Value = 'https://www.linkedin.com/in/sov-trotter/'
$ws3.Range('E7').Value = '../members/arsh.jpg'
$ws3.Range('D8').Value = 'https://www.linkedin.com/in/harshit-srivastav-1507/'
$ws3.Range('E8').Value = '../members/harshit.jpg'
$ws3.Range('D9').Value = 'https://www.linkedin.com/in/nimish-sharma-0b0929195/'
$ws3.Range('E9').Value = '../members/nimish.jpg'
$ws3.Range('D10').Value = 'https://www.linkedin.com/in/rnawathe/'
$ws3.Range('E10').Value = '../members/rohan.jpg'
$ws3.Range('D11').Value = 'https://www.linkedin.com/in/sarvesh-srivastava-03678116b/'
$ws3.Range('E11').Value = '../members/sarvesh.jpg'
$ws3.Range('D12').Value = 'https://www.facebook.com/tanuja.pal.75436'
$ws3.Range('E12').Value = '../members/tanuja.jpg'
$ws3.Range('D13').Value = 'https://www.linkedin.com/in/tanyabhandari25/'
$ws3.Range('E13').Value = '../members/tanya.jpg'
$ws3.Range('D14').Value = 'https://www.linkedin.com/in/sharma3anika/'
$ws3.Range('E14').Value = '../members/anika.jpg'
$ws3.Range('D15').Value = 'https://www.linkedin.com/in/ajay-c-200a9b110/'
$ws3.Range('E15').Value = '../members/ajay.jpg'
$ws3.Range('D16').Value = 'https://www.linkedin.com/in/vishal-dhiman-b99a9b18b/'
$ws3.Range('E16').Value = '../members/vishal.jpg'

# --- Rebuild hyperlinks on Sheet2 ---
$null = $ws2.Range('D2').Hyperlinks.Delete()
$null = $ws2.Hyperlinks.Add($ws2.Range('D2'), 'https://www.linkedin.com/in/gaurav-pandey-412180194/', '', '', 'https://www.linkedin.com/in/gaurav-pandey-412180194/')
$null = $ws2.Hyperlinks.Add($ws2.Range('D3'), 'https://www.linkedin.com/in/kunal-thakur-b27823193/', '', '', 'https://www.linkedin.com/in/kunal-thakur-b27823193/')
$null = $ws2.Hyperlinks.Add($ws2.Range('D7'), 'https://www.linkedin.com/in/mayur-kumar-47a9511a1/', '', '', 'https://www.linkedin.com/in/mayur-kumar-47a9511a1/')
$null = $ws2.Hyperlinks.Add($ws2.Range('D8'), 'https://www.linkedin.com/in/parasaggarwal/', '', '', 'https://www.linkedin.com/in/parasaggarwal/')
$null = $ws2.Hyperlinks.Add($ws2.Range('D9'), 'https://www.linkedin.com/in/parth-pant-866bb4189/', '', '', 'https://www.linkedin.com/in/parth-pant-866bb4189/')
$null = $ws2.Hyperlinks.Add($ws2.Range('D10'), 'https://github.com/ligmitz', '', '', 'https://github.com/ligmitz')
$null = $ws2.Hyperlinks.Add($ws2.Range('D13'), 'https://www.linkedin.com/in/suryansh-dwivedi-9b9975199/', '', '', 'https://www.linkedin.com/in/suryansh-dwivedi-9b9975199/')

# --- Rebuild hyperlinks on Sheet3 ---
$null = $ws3.Range('D2').Hyperlinks.Delete()
$null = $ws3.Hyperlinks.Add($ws3.Range('D4'), 'https://www.linkedin.com/in/rathod-sahaab/', '', '', 'https://www.linkedin.com/in/rathod-sahaab/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D5'), 'https://www.linkedin.com/in/aditi-singh2000/', '', '', 'https://www.linkedin.com/in/aditi-singh2000/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D7'), 'https://www.linkedin.com/in/sov-trotter/', '', '', 'https://www.linkedin.com/in/sov-trotter/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D8'), 'https://www.linkedin.com/in/harshit-srivastav-1507/', '', '', 'https://www.linkedin.com/in/harshit-srivastav-1507/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D10'), 'https://www.linkedin.com/in/rnawathe/', '', '', 'https://www.linkedin.com/in/rnawathe/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D11'), 'https://www.linkedin.com/in/sarvesh-srivastava-03678116b/', '', '', 'https://www.linkedin.com/in/sarvesh-srivastava-03678116b/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D14'), 'https://www.linkedin.com/in/sharma3anika/', '', '', 'https://www.linkedin.com/in/sharma3anika/')
$null = $ws3.Hyperlinks.Add($ws3.Range('D16'), 'https://www.linkedin.com/in/vishal-dhiman-b99a9b18b/', '', '', 'https://www.linkedin.com/in/vishal-dhiman-b99a9b18b/')

# --- Selections per sheet ---
$ws1.Activate()
$null = $ws1.Range('D17').Select()
$ws3.Activate()
$null = $ws3.Range('E17').Select()
$ws4.Activate()
$null = $ws4.Range('F23').Select()
$ws2.Activate()
$null = $ws2.Range('N18').Select()

